# Update the dSF (column F) values for specific rows as per the
# "repull data, push all data, mean calculation" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -1
    5  = 1
    10 = 8
    11 = -3
    12 = -1
    15 = -5
    17 = 12
    18 = -8
    21 = -3
    26 = -2
    27 = 6
    30 = 0
    31 = 3
    32 = 2
    34 = 3
    35 = -4
    40 = -3
    41 = 12
    42 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
